# Update the "Data" sheet visitor log:
#  - row 2 gets a new date / visitor name / approval dates
#  - two new rows (3 and 4) are appended with the same visitor info
#  - the sheet view direction is explicitly set to left-to-right
#
# All of the text-looking values (dates, names) must stay stored as TEXT
# (the workbook already keeps them as text, with numberStoredAsText
# ignored-errors), so every cell we touch gets NumberFormat "@" first to
# stop Excel from re-interpreting "08-09-2024" etc. as a real date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Explicitly set left-to-right sheet view (mirrors rightToLeft="0") ---
$excel.ActiveWindow.DisplayRightToLeft = $false

# --- Row 2: update existing visit entry ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "08-09-2024"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "hlinhbk"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2024-10-08"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2024-10-08"

# --- Row 3: new visit entry (same visitor, same day) ---
# Format the whole row as text first so the (still blank) C3:F3 cells are
# materialized in the sheet, just like C2:F2 already are in row 2.
$ws.Range("A3:H3").NumberFormat = "@"
$ws.Range("A3").Value = "08-09-2024"
$ws.Range("B3").Value = "hlinhbk"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "2024-10-08"
$ws.Range("H3").Value = "2024-10-08"

# --- Row 4: another new visit entry (same visitor, same day) ---
$ws.Range("A4:H4").NumberFormat = "@"
$ws.Range("A4").Value = "08-09-2024"
$ws.Range("B4").Value = "hlinhbk"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "2024-10-08"
$ws.Range("H4").Value = "2024-10-08"
